# Slide 14, "Use better features" / "Use a higher capacity model" content placeholder
# (5th shape on the slide): append a clarifying parenthetical in a smaller font
# to the first bullet, and grow the text box to fit the now-taller text.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(14)
$sh = $s.Shapes.Item(5)
$tf = $sh.TextFrame
$tr = $tf.TextRange

$oldFirstRun = "Use better features"
$newFirstRun = "Use better features "
$newSecondRun = "(differentiable programming to automatically learn good features)"

$para1 = $tr.Paragraphs(1)

# "Use better features" -> "Use better features " (trailing space), keeping the
# existing run's formatting (sz=2400, blue fill, +mn-lt) untouched.
$firstPart = $tr.Characters($para1.Start, $oldFirstRun.Length)
$firstPart.Text = $newFirstRun

# Append the new explanatory run right after the (now space-terminated) first run.
# Re-fetch the paragraph/range fresh since the handles above are now stale after
# the text mutation.
$tr = $sh.TextFrame.TextRange
$para1 = $tr.Paragraphs(1)
$para1.InsertAfter($newSecondRun) | Out-Null

# Restyle only the newly-inserted text to a smaller size; it inherits the color /
# typeface from the surrounding text already, so only Size needs changing.
$tr = $sh.TextFrame.TextRange
$para1 = $tr.Paragraphs(1)
$newStart = $para1.Start + $newFirstRun.Length
$newRange = $tr.Characters($newStart, $newSecondRun.Length)
$newRange.Font.Size = 16

# The shape auto-fits to its text (spAutoFit); after the edit it should end up
# cy=1144929 EMU (was 923330) while position/width stay the same.
$sh.Height = 90.1519
